# Update the Rules sheet: cell E8 ("Good Morning") becomes "GIT UPDATE",
# and the sheet's selection moves to E8 (matching the author's last
# selection when saving via jgit/commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "GIT UPDATE"
$ws.Range("E8").Select()
